# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4", before "总计".
# 2) Populate "2022-Q1" with the per-fund holding rows for that quarter.
# 3) Update the "总计" (summary) sheet: push the existing 2021-Q4 summary
#    row down to row 3 and insert a new 2022-Q1 summary row at row 2.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1) Create the new "2022-Q1" sheet right after "2021-Q4" ---------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# NOTE: look this up only *after* the insert above. Inserting a sheet shifts
# the tab position of everything after it, and a worksheet handle obtained
# beforehand would end up pointing at "2022-Q1" instead of "总计".
$total = $wb.Worksheets.Item("总计")

# Copy the header row (with its styling) and the style of the index column
# from the "2021-Q4" sheet so the new sheet matches the established format.
$q4.Range("B1:H1").Copy($q1.Range("B1:H1"))
$q4.Range("A2:A3").Copy($q1.Range("A2:A3"))

# Row 2: 华润元大景泰混合A
# Force text storage for B:G (so fund codes keep leading zeros and the
# numeric-looking figures stay strings, matching the source data), then
# drop back to the default "Normal" style so no stray number-format style
# index is left attached to the cells.
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "004976"
$q1.Range("C2").Value = "华润元大景泰混合A"
$q1.Range("D2").Value = "1.79"
$q1.Range("E2").Value = "37.61"
$q1.Range("F2").Value = "0.22"
$q1.Range("G2").Value = "0.0039"
$q1.Range("H2").Value = 8
$q1.Range("B2:G2").Style = "Normal"

# Row 3: 华润元大景泰混合C
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "004977"
$q1.Range("C3").Value = "华润元大景泰混合C"
$q1.Range("D3").Value = "1.79"
$q1.Range("E3").Value = "37.61"
$q1.Range("F3").Value = "0.22"
$q1.Range("G3").Value = "0.0039"
$q1.Range("H3").Value = 8
$q1.Range("B3:G3").Style = "Normal"

# --- 2) Update the "总计" sheet --------------------------------------------
# Push the existing 2021-Q4 summary row (row 2) down to row 3, carrying its
# styling (the bold index cell in column A) along with it.
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1

# Overwrite row 2 with the new 2022-Q1 summary figures.
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01
